$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "View" (column F, row 1) to "Cache"
$ws.Range("F1").Value = "Cache"

# Set default value of column F (F2:F26) to FALSE
$ws.Range("F2:F26").Value = $false
